# Scheduled market-data refresh for the crafting-profit sheets.
# Updates currentAveragePrice(NQ/HQ) and the derived Leve price /
# profit columns (H:N) for the rows whose backing market data moved,
# including a couple of rows where the HQ profit cell (M) needs to
# be added or removed because HQ listings appeared/disappeared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 74
$ws.Range("H74").Value = 2913.45
$ws.Range("I74").Value = 2881.4707
$ws.Range("J74").Value = 3094.6667
$ws.Range("K74").Value = 2881.4707
$ws.Range("L74").Value = 3094.6667
$ws.Range("M74").Value = -1945.4707
$ws.Range("N74").Value = -4966.6667

# Row 77
$ws.Range("H77").Value = 2913.45
$ws.Range("I77").Value = 2881.4707
$ws.Range("J77").Value = 3094.6667
$ws.Range("K77").Value = 14407.3535
$ws.Range("L77").Value = 15473.3335
$ws.Range("M77").Value = -9727.353499999999
$ws.Range("N77").Value = -24833.3335

# Row 101
$ws.Range("H101").Value = 1040.8572
$ws.Range("I101").Value = 343.2
$ws.Range("K101").Value = 1029.6
$ws.Range("M101").Value = 592.4000000000001

# Row 125
$ws.Range("H125").Value = 1935.3793
$ws.Range("J125").Value = 1776.25
$ws.Range("L125").Value = 15986.25
$ws.Range("N125").Value = -20906.25

# Row 132
$ws.Range("H132").Value = 4721419
$ws.Range("I132").Value = 5004470
$ws.Range("K132").Value = 15013410
$ws.Range("M132").Value = -15010880


$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 105127.7
$ws.Range("I45").Value = 251424.5
$ws.Range("J45").Value = 7596.5
$ws.Range("K45").Value = 251424.5
$ws.Range("L45").Value = 7596.5
$ws.Range("M45").Value = -251047.5
$ws.Range("N45").Value = -8350.5

# Row 53
$ws.Range("H53").Value = 9866.666999999999
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 9866.666999999999
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 9866.666999999999
$ws.Range("M53").ClearContents() | Out-Null
$ws.Range("N53").Value = -11230.667

# Row 74
$ws.Range("H74").Value = 1128.2307
$ws.Range("J74").Value = 1237.4
$ws.Range("L74").Value = 1237.4
$ws.Range("N74").Value = -2985.4

# Row 77
$ws.Range("H77").Value = 1128.2307
$ws.Range("J77").Value = 1237.4
$ws.Range("L77").Value = 6187
$ws.Range("N77").Value = -14923

# Row 122
$ws.Range("H122").Value = 2382.9092
$ws.Range("I122").Value = 2150.5
$ws.Range("K122").Value = 6451.5
$ws.Range("M122").Value = -4001.5


$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 1841.3191
$ws.Range("I134").Value = 1612.0444
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 4836.1332
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -2301.1332
$ws.Range("N134").Value = -26070


$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 1242.8889
$ws.Range("I5").Value = 1363.2727
$ws.Range("J5").Value = 1189.92
$ws.Range("K5").Value = 4089.8181
$ws.Range("L5").Value = 3569.76
$ws.Range("M5").Value = -3977.8181
$ws.Range("N5").Value = -3793.76

# Row 112
$ws.Range("H112").Value = 87867.914
$ws.Range("I112").Value = 334008.34
$ws.Range("J112").Value = 5821.1113
$ws.Range("K112").Value = 1002025.02
$ws.Range("L112").Value = 17463.3339
$ws.Range("M112").Value = -1000917.02
$ws.Range("N112").Value = -19679.3339

# Row 131
$ws.Range("H131").Value = 808.1900000000001
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 823.11456
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 2469.34368
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -12549.34368

# Row 135
$ws.Range("H135").Value = 1242.8889
$ws.Range("I135").Value = 1363.2727
$ws.Range("J135").Value = 1189.92
$ws.Range("K135").Value = 12269.4543
$ws.Range("L135").Value = 10709.28
$ws.Range("M135").Value = -9734.454299999999
$ws.Range("N135").Value = -15779.28


$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 46840.28
$ws.Range("I70").Value = 56700.55
$ws.Range("J70").Value = 7399.2
$ws.Range("K70").Value = 56700.55
$ws.Range("L70").Value = 7399.2
$ws.Range("M70").Value = -56430.55
$ws.Range("N70").Value = -7939.2

# Row 73
$ws.Range("H73").Value = 46840.28
$ws.Range("I73").Value = 56700.55
$ws.Range("J73").Value = 7399.2
$ws.Range("K73").Value = 56700.55
$ws.Range("L73").Value = 7399.2
$ws.Range("M73").Value = -55764.55
$ws.Range("N73").Value = -9271.200000000001

# Row 102
$ws.Range("H102").Value = 3537.182
$ws.Range("I102").Value = 3676
$ws.Range("K102").Value = 3676
$ws.Range("M102").Value = -2054

# Row 126
$ws.Range("H126").Value = 2647.2942
$ws.Range("I126").Value = 3054.5557
$ws.Range("J126").Value = 2189.125
$ws.Range("K126").Value = 9163.667099999999
$ws.Range("L126").Value = 6567.375
$ws.Range("M126").Value = -6693.667099999999
$ws.Range("N126").Value = -11507.375


$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 70358.87
$ws.Range("I40").Value = 172300.5
$ws.Range("J40").Value = 2397.7778
$ws.Range("K40").Value = 172300.5
$ws.Range("L40").Value = 2397.7778
$ws.Range("M40").Value = -172164.5
$ws.Range("N40").Value = -2669.7778

# Row 68
$ws.Range("H68").Value = 4247.6924
$ws.Range("I68").Value = 1725
$ws.Range("J68").Value = 5368.8887
$ws.Range("K68").Value = 1725
$ws.Range("L68").Value = 5368.8887
$ws.Range("M68").Value = -976
$ws.Range("N68").Value = -6866.8887

# Row 71
$ws.Range("H71").Value = 4247.6924
$ws.Range("I71").Value = 1725
$ws.Range("J71").Value = 5368.8887
$ws.Range("K71").Value = 8625
$ws.Range("L71").Value = 26844.4435
$ws.Range("M71").Value = -4881
$ws.Range("N71").Value = -34332.4435

# Row 82
$ws.Range("H82").Value = 1979.6666
$ws.Range("I82").Value = 980
$ws.Range("J82").Value = 2070.5454
$ws.Range("K82").Value = 980
$ws.Range("L82").Value = 2070.5454
$ws.Range("N82").Value = -2792.5454
$ws.Range("M82").Value = -619

# Row 85
$ws.Range("H85").Value = 1979.6666
$ws.Range("I85").Value = 980
$ws.Range("J85").Value = 2070.5454
$ws.Range("K85").Value = 980
$ws.Range("L85").Value = 2070.5454
$ws.Range("N85").Value = -4566.5454
$ws.Range("M85").Value = 268

# Row 93
$ws.Range("H93").Value = 1395.9
$ws.Range("I93").Value = 1277.7241
$ws.Range("J93").Value = 1707.4546
$ws.Range("K93").Value = 1277.7241
$ws.Range("L93").Value = 1707.4546
$ws.Range("M93").Value = -29.72409999999991
$ws.Range("N93").Value = -4203.4546

# Row 100
$ws.Range("H100").Value = 1788.4286
$ws.Range("I100").Value = 1515
$ws.Range("J100").Value = 1993.5
$ws.Range("K100").Value = 1515
$ws.Range("L100").Value = 1993.5
$ws.Range("M100").Value = -974
$ws.Range("N100").Value = -3075.5


$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Range("H107").Value = 999999.5
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents() | Out-Null

# Row 122
$ws.Range("H122").Value = 2726.3125
$ws.Range("I122").Value = 2250
$ws.Range("J122").Value = 2794.3572
$ws.Range("K122").Value = 6750
$ws.Range("L122").Value = 8383.071599999999
$ws.Range("M122").Value = -4300
$ws.Range("N122").Value = -13283.0716

# Row 126
$ws.Range("H126").Value = 1752.1111
$ws.Range("I126").Value = 1516
$ws.Range("J126").Value = 2047.25
$ws.Range("K126").Value = 4548
$ws.Range("L126").Value = 6141.75
$ws.Range("M126").Value = -2078
$ws.Range("N126").Value = -11081.75

# Row 132
$ws.Range("H132").Value = 3001.8
$ws.Range("I132").Value = 4121.0527
$ws.Range("J132").Value = 1672.6875
$ws.Range("K132").Value = 12363.1581
$ws.Range("L132").Value = 5018.0625
$ws.Range("M132").Value = -9833.158100000001
$ws.Range("N132").Value = -10078.0625
